$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

for ($r = 8; $r -le 13; $r++) {
    # Priority column (E) on the zh-cn and de-de sheets: "" -> "ht"
    $zh.Cells.Item($r, 5).Value = "ht"
    $de.Cells.Item($r, 5).Value = "ht"

    # Latest Handoff Datetime column (H) on the zh-cn sheet
    $zh.Cells.Item($r, 8).Value = "2016-09-08 04:30:49"

    # Latest Handoff Datetime column (H) on the de-de sheet, and the
    # corresponding Latest HO Xliff Generate Date column (G) on the
    # Overview sheet (backed by the same shared string in the source report)
    $de.Cells.Item($r, 8).Value = "2016-09-08 04:30:55"
    $overview.Cells.Item($r, 7).Value = "2016-09-08 04:30:55"
}
